$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 448.55225
$ws.Range("J17").Value = 305.98306
$ws.Range("L17").Value = 917.9491800000001
$ws.Range("N17").Value = -1253.94918

$ws.Range("H33").Value = 246.45454
$ws.Range("I33").Value = 242.75
$ws.Range("J33").Value = 256.33334
$ws.Range("K33").Value = 242.75
$ws.Range("L33").Value = 256.33334
$ws.Range("M33").Value = -13.75
$ws.Range("N33").Value = -714.33334

$ws.Range("H62").Value = 9886
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 9886
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 9886
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -11134

$ws.Range("H65").Value = 9886
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 9886
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 49430
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -55670

$ws.Range("H88").Value = 3600
$ws.Range("I88").Value = 5900
$ws.Range("J88").Value = 2450
$ws.Range("K88").Value = 5900
$ws.Range("L88").Value = 2450
$ws.Range("M88").Value = -5494
$ws.Range("N88").Value = -3262

$ws.Range("H91").Value = 3600
$ws.Range("I91").Value = 5900
$ws.Range("J91").Value = 2450
$ws.Range("K91").Value = 5900
$ws.Range("L91").Value = 2450
$ws.Range("M91").Value = -4496
$ws.Range("N91").Value = -5258

$ws.Range("H100").Value = 28573506
$ws.Range("I100").Value = 28573506
$ws.Range("K100").Value = 28573506
$ws.Range("M100").Value = -28572965

$ws.Range("H106").Value = 45459210
$ws.Range("I106").Value = 2697
$ws.Range("K106").Value = 2697
$ws.Range("M106").Value = -2066

$ws.Range("H137").Value = 1537980.1
$ws.Range("I137").Value = 1985136.5
$ws.Range("J137").Value = 4872.5713
$ws.Range("K137").Value = 5955409.5
$ws.Range("L137").Value = 14617.7139
$ws.Range("M137").Value = -5952859.5
$ws.Range("N137").Value = -19717.7139

$ws.Range("H138").Value = 2844.14
$ws.Range("I138").Value = 1442.5
$ws.Range("J138").Value = 2902.5417
$ws.Range("K138").Value = 4327.5
$ws.Range("L138").Value = 8707.625100000001
$ws.Range("M138").Value = 812.5
$ws.Range("N138").Value = -18987.6251

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 981.25
$ws.Range("I2").Value = 1456
$ws.Range("J2").Value = 506.5
$ws.Range("K2").Value = 1456
$ws.Range("L2").Value = 506.5
$ws.Range("M2").Value = -1343
$ws.Range("N2").Value = -732.5

$ws.Range("H45").Value = 2595
$ws.Range("I45").Value = 2380
$ws.Range("J45").Value = 2666.6667
$ws.Range("K45").Value = 2380
$ws.Range("L45").Value = 2666.6667
$ws.Range("M45").Value = -2003
$ws.Range("N45").Value = -3420.6667

$ws.Range("H116").Value = 981.25
$ws.Range("I116").Value = 1456
$ws.Range("J116").Value = 506.5
$ws.Range("K116").Value = 1456
$ws.Range("L116").Value = 506.5
$ws.Range("M116").Value = 838
$ws.Range("N116").Value = -5094.5

$ws.Range("H122").Value = 4182.4
$ws.Range("I122").Value = 1304
$ws.Range("J122").Value = 8500
$ws.Range("K122").Value = 3912
$ws.Range("L122").Value = 25500
$ws.Range("M122").Value = -1462
$ws.Range("N122").Value = -30400

$ws.Range("H137").Value = 43378.332
$ws.Range("J137").Value = 43378.332
$ws.Range("L137").Value = 43378.332
$ws.Range("N137").Value = -53578.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 981.25
$ws.Range("I3").Value = 1456
$ws.Range("J3").Value = 506.5
$ws.Range("K3").Value = 1456
$ws.Range("L3").Value = 506.5
$ws.Range("M3").Value = -1342
$ws.Range("N3").Value = -734.5

$ws.Range("H64").Value = 817.9545000000001
$ws.Range("J64").Value = 836.8461
$ws.Range("L64").Value = 836.8461
$ws.Range("N64").Value = -1286.8461

$ws.Range("H67").Value = 817.9545000000001
$ws.Range("J67").Value = 836.8461
$ws.Range("L67").Value = 836.8461
$ws.Range("N67").Value = -2396.8461

$ws.Range("H80").Value = 660.3333
$ws.Range("I80").Value = 1055
$ws.Range("J80").Value = 581.4
$ws.Range("K80").Value = 1055
$ws.Range("L80").Value = 581.4
$ws.Range("M80").Value = -57
$ws.Range("N80").Value = -2577.4

$ws.Range("H83").Value = 660.3333
$ws.Range("I83").Value = 1055
$ws.Range("J83").Value = 581.4
$ws.Range("K83").Value = 5275
$ws.Range("L83").Value = 2907
$ws.Range("M83").Value = -283
$ws.Range("N83").Value = -12891

$ws.Range("H134").Value = 1876.9459
$ws.Range("I134").Value = 1156.7931
$ws.Range("K134").Value = 3470.379300000001
$ws.Range("M134").Value = -935.3793000000005

$ws.Range("H137").Value = 45310
$ws.Range("J137").Value = 45310
$ws.Range("L137").Value = 45310
$ws.Range("N137").Value = -55510

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4795.8
$ws.Range("I31").Value = 1480.75
$ws.Range("K31").Value = 1480.75
$ws.Range("M31").Value = -1185.75

$ws.Range("H34").Value = 4795.8
$ws.Range("I34").Value = 1480.75
$ws.Range("K34").Value = 1480.75
$ws.Range("M34").Value = -1278.75

$ws.Range("H58").Value = 3006.0334
$ws.Range("I58").Value = 1828.1459
$ws.Range("J58").Value = 7717.5835
$ws.Range("K58").Value = 1828.1459
$ws.Range("L58").Value = 7717.5835
$ws.Range("M58").Value = -1625.1459
$ws.Range("N58").Value = -8123.5835

$ws.Range("H86").Value = 2307.25
$ws.Range("I86").Value = 1989.2
$ws.Range("J86").Value = 3261.4
$ws.Range("K86").Value = 1989.2
$ws.Range("L86").Value = 3261.4
$ws.Range("M86").Value = -866.2
$ws.Range("N86").Value = -5507.4

$ws.Range("H89").Value = 2307.25
$ws.Range("I89").Value = 1989.2
$ws.Range("J89").Value = 3261.4
$ws.Range("K89").Value = 9946
$ws.Range("L89").Value = 16307
$ws.Range("M89").Value = -4330
$ws.Range("N89").Value = -27539

$ws.Range("H122").Value = 2783.4375
$ws.Range("I122").Value = 2042.3
$ws.Range("J122").Value = 4018.6667
$ws.Range("K122").Value = 6126.9
$ws.Range("L122").Value = 12056.0001
$ws.Range("M122").Value = -3676.9
$ws.Range("N122").Value = -16956.0001

$ws.Range("H136").Value = 3006.0334
$ws.Range("I136").Value = 1828.1459
$ws.Range("J136").Value = 7717.5835
$ws.Range("K136").Value = 5484.4377
$ws.Range("L136").Value = 23152.7505
$ws.Range("M136").Value = -2934.4377
$ws.Range("N136").Value = -28252.7505

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 3000
$ws.Range("I47").Value = 1000
$ws.Range("J47").Value = 5000
$ws.Range("K47").Value = 3000
$ws.Range("L47").Value = 15000
$ws.Range("M47").Value = -2569
$ws.Range("N47").Value = -15862

$ws.Range("I113").Value = 632.2353000000001
$ws.Range("J113").Value = 8333934.5
$ws.Range("K113").Value = 1896.7059
$ws.Range("L113").Value = 25001803.5
$ws.Range("M113").Value = 273.2940999999998
$ws.Range("N113").Value = -25006143.5

$ws.Range("H122").Value = 2593.432
$ws.Range("I122").Value = 800.1875
$ws.Range("J122").Value = 3618.1428
$ws.Range("K122").Value = 7201.6875
$ws.Range("L122").Value = 32563.2852
$ws.Range("M122").Value = -4751.6875
$ws.Range("N122").Value = -37463.2852

$ws.Range("H129").Value = 3218.111
$ws.Range("I129").Value = 5207.5
$ws.Range("J129").Value = 1626.6
$ws.Range("K129").Value = 15622.5
$ws.Range("L129").Value = 4879.799999999999
$ws.Range("M129").Value = -10622.5
$ws.Range("N129").Value = -14879.8

$ws.Range("H131").Value = 769.36
$ws.Range("I131").Value = 355
$ws.Range("J131").Value = 800.5484
$ws.Range("K131").Value = 1065
$ws.Range("L131").Value = 2401.6452
$ws.Range("M131").Value = 3975
$ws.Range("N131").Value = -12481.6452

$ws.Range("H138").Value = 3500
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3500
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 10500
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -20780

$ws.Range("H139").Value = 1630.25
$ws.Range("I139").Value = 1013.6667
$ws.Range("J139").Value = 3480
$ws.Range("K139").Value = 3041.0001
$ws.Range("L139").Value = 10440
$ws.Range("M139").Value = 2098.9999
$ws.Range("N139").Value = -20720

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 23533.334
$ws.Range("J45").Value = 23533.334
$ws.Range("L45").Value = 23533.334
$ws.Range("N45").Value = -24651.334

$ws.Range("H46").Value = 33164.668
$ws.Range("J46").Value = 33164.668
$ws.Range("L46").Value = 33164.668
$ws.Range("N46").Value = -33476.668

$ws.Range("H122").Value = 5748.1
$ws.Range("I122").Value = 2745.5
$ws.Range("J122").Value = 10252
$ws.Range("K122").Value = 8236.5
$ws.Range("L122").Value = 30756
$ws.Range("M122").Value = -5786.5
$ws.Range("N122").Value = -35656

$ws.Range("H137").Value = 37730
$ws.Range("J137").Value = 37730
$ws.Range("L137").Value = 37730
$ws.Range("N137").Value = -47930

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7547.75
$ws.Range("I40").Value = 3474.6667
$ws.Range("K40").Value = 3474.6667
$ws.Range("M40").Value = -3338.6667

$ws.Range("H61").Value = 2146.5334
$ws.Range("I61").Value = 2534
$ws.Range("J61").Value = 1807.5
$ws.Range("K61").Value = 2534
$ws.Range("L61").Value = 1807.5
$ws.Range("M61").Value = -2332
$ws.Range("N61").Value = -2211.5

$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490

$ws.Range("H113").Value = 2146.5334
$ws.Range("I113").Value = 2534
$ws.Range("J113").Value = 1807.5
$ws.Range("K113").Value = 2534
$ws.Range("L113").Value = 1807.5
$ws.Range("M113").Value = -364
$ws.Range("N113").Value = -6147.5

$ws.Range("H122").Value = 3750.5652
$ws.Range("I122").Value = 2098.625
$ws.Range("J122").Value = 4631.6
$ws.Range("K122").Value = 6295.875
$ws.Range("L122").Value = 13894.8
$ws.Range("M122").Value = -3845.875
$ws.Range("N122").Value = -18794.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 34997.645
$ws.Range("J103").Value = 34997.645
$ws.Range("L103").Value = 34997.645
$ws.Range("N103").Value = -37341.645

$ws.Range("H107").Value = 175.5
$ws.Range("I107").Value = 162.5
$ws.Range("J107").Value = 201.5
$ws.Range("K107").Value = 487.5
$ws.Range("L107").Value = 604.5
$ws.Range("M107").Value = 1432.5
$ws.Range("N107").Value = -4444.5
